$d = $word.ActiveDocument

$replacements = @(
    @("197×9=1773", "322×2=644"),
    @("948×7=6636", "537×2=1074"),
    @("227×2=454", "518×6=3108"),
    @("371×7=2597", "276×9=2484"),
    @("917×4=3668", "871×3=2613"),
    @("915×6=5490", "638×2=1276"),
    @("230×8=1840", "964×2=1928"),
    @("682×8=5456", "120×3=360"),
    @("383×5=1915", "373×8=2984"),
    @("232×4=928", "700×8=5600"),
    @("264×8=2112", "460×6=2760"),
    @("451×8=3608", "740×9=6660"),
    @("342×6=2052", "538×9=4842"),
    @("857×2=1714", "846×7=5922"),
    @("484×6=2904", "237×3=711"),
    @("405×2=810", "316×5=1580"),
    @("255×3=765", "551×5=2755"),
    @("451×3=1353", "236×2=472"),
    @("272×6=1632", "893×6=5358"),
    @("223×9=2007", "754×8=6032"),
    @("313×9=2817", "815×5=4075"),
    @("974×2=1948", "603×2=1206"),
    @("624×7=4368", "266×8=2128"),
    @("540×9=4860", "989×7=6923"),
    @("933×8=7464", "199×2=398")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
